# Auto-generated edit script applying market price / profit recalculations
# to the Leve profit tables across all job sheets, per the scheduled runner update.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4260
$ws.Range("I64").Value = 4260
$ws.Range("K64").Value = 4260
$ws.Range("M64").Value = -4012
$ws.Range("H67").Value = 4260
$ws.Range("I67").Value = 4260
$ws.Range("K67").Value = 4260
$ws.Range("M67").Value = -3402
$ws.Range("H74").Value = 4043.75
$ws.Range("I74").Value = 4133.3335
$ws.Range("J74").Value = 3775
$ws.Range("K74").Value = 4133.3335
$ws.Range("L74").Value = 3775
$ws.Range("M74").Value = -3197.3335
$ws.Range("N74").Value = -5647
$ws.Range("H76").Value = 3365
$ws.Range("I76").Value = 3542.8572
$ws.Range("J76").Value = 3116
$ws.Range("K76").Value = 3542.8572
$ws.Range("L76").Value = 3116
$ws.Range("M76").Value = -3227.8572
$ws.Range("N76").Value = -3746
$ws.Range("H77").Value = 4043.75
$ws.Range("I77").Value = 4133.3335
$ws.Range("J77").Value = 3775
$ws.Range("K77").Value = 20666.6675
$ws.Range("L77").Value = 18875
$ws.Range("M77").Value = -15986.6675
$ws.Range("N77").Value = -28235
$ws.Range("H79").Value = 3365
$ws.Range("I79").Value = 3542.8572
$ws.Range("J79").Value = 3116
$ws.Range("K79").Value = 3542.8572
$ws.Range("L79").Value = 3116
$ws.Range("M79").Value = -2450.8572
$ws.Range("N79").Value = -5300
$ws.Range("H92").Value = 655.0454999999999
$ws.Range("I92").Value = 541.0625
$ws.Range("J92").Value = 959
$ws.Range("K92").Value = 541.0625
$ws.Range("L92").Value = 959
$ws.Range("M92").Value = 706.9375
$ws.Range("N92").Value = -3455
$ws.Range("H132").Value = 4773.8096
$ws.Range("I132").Value = 4773.8096
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 14321.4288
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -11791.4288
$ws.Range("N132").ClearContents()
$ws.Range("H137").Value = 29415286
$ws.Range("I137").Value = 1936.0714
$ws.Range("J137").Value = 50004630
$ws.Range("K137").Value = 5808.2142
$ws.Range("L137").Value = 150013890
$ws.Range("M137").Value = -3258.2142
$ws.Range("N137").Value = -150018990

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 21999
$ws.Range("J23").Value = 21999
$ws.Range("L23").Value = 21999
$ws.Range("N23").Value = -22517
$ws.Range("H32").Value = 17412.605
$ws.Range("I32").Value = 14675.594
$ws.Range("J32").Value = 104997
$ws.Range("K32").Value = 14675.594
$ws.Range("L32").Value = 104997
$ws.Range("M32").Value = -14388.594
$ws.Range("N32").Value = -105571
$ws.Range("H37").Value = 10123.777
$ws.Range("J37").Value = 10123.777
$ws.Range("L37").Value = 10123.777
$ws.Range("N37").Value = -10669.777
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H61").Value = 5293757.5
$ws.Range("I61").Value = 6538183
$ws.Range("K61").Value = 6538183
$ws.Range("M61").Value = -6537971
$ws.Range("H63").Value = 6281.1113
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 6281.1113
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 6281.1113
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -7653.1113
$ws.Range("H66").Value = 6281.1113
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 6281.1113
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 31405.5565
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -38269.5565
$ws.Range("H132").Value = 2466804.5
$ws.Range("I132").Value = 4314486
$ws.Range("J132").Value = 3229.3333
$ws.Range("K132").Value = 12943458
$ws.Range("L132").Value = 9687.999899999999
$ws.Range("M132").Value = -12940928
$ws.Range("N132").Value = -14747.9999
$ws.Range("H136").Value = 5293757.5
$ws.Range("I136").Value = 6538183
$ws.Range("K136").Value = 19614549
$ws.Range("M136").Value = -19611999
$ws.Range("H138").Value = 66397.336
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 66397.336
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 66397.336
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -76677.336

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2377.75
$ws.Range("I105").Value = 2100
$ws.Range("J105").Value = 2655.5
$ws.Range("K105").Value = 2100
$ws.Range("L105").Value = 2655.5
$ws.Range("M105").Value = -353
$ws.Range("N105").Value = -6149.5
$ws.Range("H134").Value = 28575442
$ws.Range("I134").Value = 28575442
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 85726326
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -85723791
$ws.Range("N134").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6085.294
$ws.Range("I31").Value = 2110.3
$ws.Range("J31").Value = 7741.5415
$ws.Range("K31").Value = 2110.3
$ws.Range("L31").Value = 7741.5415
$ws.Range("M31").Value = -1815.3
$ws.Range("N31").Value = -8331.541499999999
$ws.Range("H34").Value = 6085.294
$ws.Range("I34").Value = 2110.3
$ws.Range("J34").Value = 7741.5415
$ws.Range("K34").Value = 2110.3
$ws.Range("L34").Value = 7741.5415
$ws.Range("M34").Value = -1908.3
$ws.Range("N34").Value = -8145.5415
$ws.Range("H134").Value = 2294.0977
$ws.Range("I134").Value = 1919.3636
$ws.Range("J134").Value = 3839.875
$ws.Range("K134").Value = 5758.0908
$ws.Range("L134").Value = 11519.625
$ws.Range("M134").Value = -3223.0908
$ws.Range("N134").Value = -16589.625
$ws.Range("H141").Value = 57029.637
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 57029.637
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 57029.637
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -67389.637

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 50730080
$ws.Range("J131").Value = 26319218
$ws.Range("L131").Value = 78957654
$ws.Range("N131").Value = -78967734

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4890.769
$ws.Range("I70").Value = 5316
$ws.Range("J70").Value = 4625
$ws.Range("K70").Value = 5316
$ws.Range("L70").Value = 4625
$ws.Range("M70").Value = -5046
$ws.Range("N70").Value = -5165
$ws.Range("H73").Value = 4890.769
$ws.Range("I73").Value = 5316
$ws.Range("J73").Value = 4625
$ws.Range("K73").Value = 5316
$ws.Range("L73").Value = 4625
$ws.Range("M73").Value = -4380
$ws.Range("N73").Value = -6497
$ws.Range("H80").Value = 140568.25
$ws.Range("I80").Value = 4113.3335
$ws.Range("J80").Value = 222441.2
$ws.Range("K80").Value = 4113.3335
$ws.Range("L80").Value = 222441.2
$ws.Range("M80").Value = -3115.3335
$ws.Range("N80").Value = -224437.2
$ws.Range("H83").Value = 140568.25
$ws.Range("I83").Value = 4113.3335
$ws.Range("J83").Value = 222441.2
$ws.Range("K83").Value = 20566.6675
$ws.Range("L83").Value = 1112206
$ws.Range("M83").Value = -15574.6675
$ws.Range("N83").Value = -1122190
$ws.Range("H132").Value = 2177.75
$ws.Range("I132").Value = 1457.8
$ws.Range("J132").Value = 3377.6667
$ws.Range("K132").Value = 4373.4
$ws.Range("L132").Value = 10133.0001
$ws.Range("M132").Value = -1843.4
$ws.Range("N132").Value = -15193.0001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2166.6924
$ws.Range("I7").Value = 2206.182
$ws.Range("J7").Value = 1949.5
$ws.Range("K7").Value = 2206.182
$ws.Range("L7").Value = 1949.5
$ws.Range("M7").Value = -2094.182
$ws.Range("N7").Value = -2173.5
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H68").Value = 2124.138
$ws.Range("I68").Value = 1902.9412
$ws.Range("K68").Value = 1902.9412
$ws.Range("M68").Value = -1153.9412
$ws.Range("H70").Value = 29666.666
$ws.Range("J70").Value = 29500
$ws.Range("L70").Value = 29500
$ws.Range("N70").Value = -30040
$ws.Range("H71").Value = 2124.138
$ws.Range("I71").Value = 1902.9412
$ws.Range("K71").Value = 9514.706
$ws.Range("M71").Value = -5770.706
$ws.Range("H73").Value = 29666.666
$ws.Range("J73").Value = 29500
$ws.Range("L73").Value = 29500
$ws.Range("N73").Value = -31372
$ws.Range("H122").Value = 1967.3489
$ws.Range("I122").Value = 1841.6562
$ws.Range("J122").Value = 2333
$ws.Range("K122").Value = 5524.9686
$ws.Range("L122").Value = 6999
$ws.Range("M122").Value = -3074.9686
$ws.Range("N122").Value = -11899
$ws.Range("H126").Value = 2166.6924
$ws.Range("I126").Value = 2206.182
$ws.Range("J126").Value = 1949.5
$ws.Range("K126").Value = 6618.545999999999
$ws.Range("L126").Value = 5848.5
$ws.Range("M126").Value = -4148.545999999999
$ws.Range("N126").Value = -10788.5
$ws.Range("H130").Value = 49139.668
$ws.Range("J130").Value = 49139.668
$ws.Range("L130").Value = 49139.668
$ws.Range("N130").Value = -59179.668
$ws.Range("H137").Value = 49726
$ws.Range("I137").Value = 35195
$ws.Range("J137").Value = 59413.332
$ws.Range("K137").Value = 35195
$ws.Range("L137").Value = 59413.332
$ws.Range("M137").Value = -30095
$ws.Range("N137").Value = -69613.33199999999
